# Revert to old version of dispatch priority
$wb = $excel.ActiveWorkbook

$wsAbout  = $wb.Worksheets.Item("About")
$wsBDPbES = $wb.Worksheets.Item("BDPbES")

# --- About sheet: drop the "New Mexico" label cell and roll back the date ---
$wsAbout.Range("B1").ClearContents()
$wsAbout.Range("C1").Value = 44307

# --- BDPbES sheet: restore old priority values for onshore wind & solar PV ---
$wsBDPbES.Range("B6").Value = 2
$wsBDPbES.Range("B7").Value = 2

# --- Restore the active/selected tab to "About" (was "BDPbES") ---
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
